$d = $word.ActiveDocument

# Locate the "graded 200+ assignments" bullet and bump the grading count
# from 200 to 400 (the "+" and the following text stay as-is).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# MatchWholeWord = $true keeps this from touching any other digit
# sequence (e.g. rsids, years) that merely contains "200".
$find.Execute("200", $true, $true, $false, $false, $false, $true, 1, $false, "400", 2)
